$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): replace date headers with homework labels ---
$ws.Range("C2").Value = "ДЗ_1"
$ws.Range("D2").Value = "ДЗ_2"
$ws.Range("E2").Value = "ДЗ_3"
$ws.Range("F2").Value = "ДЗ_4"
$ws.Range("G2").Value = "ДЗ_5"
$ws.Range("H2").ClearContents()

# --- Clear all the attendance/grade marks in the data rows (4-30), columns C:H ---
$ws.Range("C4:H30").ClearContents()

# --- Clear the totals row (31) ---
$ws.Range("C31:H31").ClearContents()

# --- View state: freeze panes still split at column B / row 3, but scrolled so
#     row 18 is the first visible row under the frozen header; selection is C2:G2 ---
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("C18").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("C4").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C2:G2").Select()
